$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 372, pushing existing rows 372..479 down to 373..480
$ws.Rows.Item(372).Insert()

# Populate the newly inserted row 372 with the new record
$ws.Range("A372").Value = 11
$ws.Range("B372").Value = 'Vega Monumental Concepción'
$ws.Range("C372").Value = 'Bíobío'
$ws.Range("D372").Value = 44663
$ws.Range("E372").Value = 8
$ws.Range("F372").Value = 'Fruta'
$ws.Range("G372").Value = 100108
$ws.Range("H372").Value = 'Tropicales y subtropicales'
$ws.Range("I372").Value = 100108006
$ws.Range("J372").Value = 'Plátano'
$ws.Range("K372").Value = 'Sin especificar'
$ws.Range("L372").Value = 'Pintón'
$ws.Range("M372").Value = 700
$ws.Range("N372").Value = 15000
$ws.Range("O372").Value = 16000
$ws.Range("P372").Value = 15571
$ws.Range("Q372").Value = '$/caja 20 kilos'
$ws.Range("R372").Value = 'Ecuador'
$ws.Range("S372").Value = 779
$ws.Range("T372").Value = 20
